$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row 8: period labels ---
$ws.Range("D8").Value = "12 ماهه منتهی به 1397/12"
$ws.Range("E8").Value = "12 ماهه منتهی به 1398/12"
$ws.Range("F8").Value = "12 ماهه منتهی به 1399/12"
$ws.Range("G8").Value = "12 ماهه منتهی به 1400/12"
$ws.Range("H8").Value = "12 ماهه منتهی به 1401/12"

# --- Header row 9: publish dates ---
$ws.Range("D9").Value = "1399-04-18 (8)"
$ws.Range("E9").Value = "1400-04-19 (8)"
$ws.Range("F9").Value = "1401-05-19 (9)"
$ws.Range("G9").Value = "1402-02-27 (7)"
$ws.Range("H9").Value = "1402-02-27"

# --- Data rows 12-58 ---
# Row 12
$ws.Range("D12").Value = 57264
$ws.Range("E12").Value = 247513
$ws.Range("F12").Value = 301092
$ws.Range("G12").Value = 443348
$ws.Range("H12").Value = 1908937
# Row 13
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 0
$ws.Range("H13").Value = 0
# Row 14
$ws.Range("D14").Value = 600377
$ws.Range("E14").Value = 369775
$ws.Range("F14").Value = 228175
$ws.Range("G14").Value = 2018602
$ws.Range("H14").Value = 2479032
# Row 15
$ws.Range("D15").Value = 351273
$ws.Range("E15").Value = 542524
$ws.Range("F15").Value = 856049
$ws.Range("G15").Value = 2364123
$ws.Range("H15").Value = 2752769
# Row 16
$ws.Range("D16").Value = 221205
$ws.Range("E16").Value = 85679
$ws.Range("F16").Value = 349698
$ws.Range("G16").Value = 158076
$ws.Range("H16").Value = 350174
# Row 17
$ws.Range("D17").Value = 0
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 0
$ws.Range("G17").Value = 0
$ws.Range("H17").Value = 0
# Row 18
$ws.Range("D18").Value = 1230119
$ws.Range("E18").Value = 1245491
$ws.Range("F18").Value = 1735014
$ws.Range("G18").Value = 4984149
$ws.Range("H18").Value = 7490912
# Row 19
$ws.Range("D19").Value = 0
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 0
$ws.Range("G19").Value = 0
$ws.Range("H19").Value = 0
# Row 20
$ws.Range("D20").Value = 0
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 0
$ws.Range("G20").Value = 0
$ws.Range("H20").Value = 0
# Row 21
$ws.Range("D21").Value = 0
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 0
$ws.Range("G21").Value = 0
$ws.Range("H21").Value = 0
# Row 22
$ws.Range("D22").Value = 126628
$ws.Range("E22").Value = 276909
$ws.Range("F22").Value = 1220254
$ws.Range("G22").Value = 1597772
$ws.Range("H22").Value = 3010651
# Row 23
$ws.Range("D23").Value = 12354
$ws.Range("E23").Value = 12282
$ws.Range("F23").Value = 12108
$ws.Range("G23").Value = 11113
$ws.Range("H23").Value = 12463
# Row 24
$ws.Range("D24").Value = "-"
$ws.Range("E24").Value = "-"
$ws.Range("F24").Value = "-"
$ws.Range("G24").Value = "-"
$ws.Range("H24").Value = "-"
# Row 25
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 0
$ws.Range("G25").Value = 0
$ws.Range("H25").Value = 12500
# Row 26
$ws.Range("D26").Value = 138982
$ws.Range("E26").Value = 289191
$ws.Range("F26").Value = 1232362
$ws.Range("G26").Value = 1608885
$ws.Range("H26").Value = 3035614
# Row 27
$ws.Range("D27").Value = 1369101
$ws.Range("E27").Value = 1534682
$ws.Range("F27").Value = 2967376
$ws.Range("G27").Value = 6593034
$ws.Range("H27").Value = 10526526
# Row 29
$ws.Range("D29").Value = 62243
$ws.Range("E29").Value = 196222
$ws.Range("F29").Value = 890995
$ws.Range("G29").Value = 1944356
$ws.Range("H29").Value = 4374903
# Row 30
$ws.Range("D30").Value = "-"
$ws.Range("E30").Value = "-"
$ws.Range("F30").Value = "-"
$ws.Range("G30").Value = "-"
$ws.Range("H30").Value = "-"
# Row 31
$ws.Range("D31").Value = 199370
$ws.Range("E31").Value = 164690
$ws.Range("F31").Value = 411185
$ws.Range("G31").Value = 80778
$ws.Range("H31").Value = 587319
# Row 32
$ws.Range("D32").Value = 835
$ws.Range("E32").Value = 3297
$ws.Range("F32").Value = 1496
$ws.Range("G32").Value = 1496
$ws.Range("H32").Value = 30
# Row 33
$ws.Range("D33").Value = 0
$ws.Range("E33").Value = 272782
$ws.Range("F33").Value = 208881
$ws.Range("G33").Value = 532053
$ws.Range("H33").Value = 178104
# Row 34
$ws.Range("D34").Value = 33903
$ws.Range("E34").Value = 48001
$ws.Range("F34").Value = 33702
$ws.Range("G34").Value = 865372
$ws.Range("H34").Value = 1006569
# Row 35
$ws.Range("D35").Value = 0
$ws.Range("E35").Value = 0
$ws.Range("F35").Value = 0
$ws.Range("G35").Value = 0
$ws.Range("H35").Value = 0
# Row 36
$ws.Range("D36").Value = 0
$ws.Range("E36").Value = 0
$ws.Range("F36").Value = 0
$ws.Range("G36").Value = 0
$ws.Range("H36").Value = 0
# Row 37
$ws.Range("D37").Value = 296351
$ws.Range("E37").Value = 684992
$ws.Range("F37").Value = 1546259
$ws.Range("G37").Value = 3424055
$ws.Range("H37").Value = 6146925
# Row 38
$ws.Range("D38").Value = 0
$ws.Range("E38").Value = 0
$ws.Range("F38").Value = 0
$ws.Range("G38").Value = 0
$ws.Range("H38").Value = 0
# Row 39
$ws.Range("D39").Value = "-"
$ws.Range("E39").Value = "-"
$ws.Range("F39").Value = "-"
$ws.Range("G39").Value = "-"
$ws.Range("H39").Value = "-"
# Row 40
$ws.Range("D40").Value = 71153
$ws.Range("E40").Value = 52082
$ws.Range("F40").Value = 18380
$ws.Range("G40").Value = 0
$ws.Range("H40").Value = 0
# Row 41
$ws.Range("D41").Value = 12645
$ws.Range("E41").Value = 18982
$ws.Range("F41").Value = 29770
$ws.Range("G41").Value = 47476
$ws.Range("H41").Value = 86802
# Row 42
$ws.Range("D42").Value = 83798
$ws.Range("E42").Value = 71064
$ws.Range("F42").Value = 48150
$ws.Range("G42").Value = 47476
$ws.Range("H42").Value = 86802
# Row 43
$ws.Range("D43").Value = 380149
$ws.Range("E43").Value = 756056
$ws.Range("F43").Value = 1594409
$ws.Range("G43").Value = 3471531
$ws.Range("H43").Value = 6233727
# Row 45
$ws.Range("D45").Value = 200000
$ws.Range("E45").Value = 200000
$ws.Range("F45").Value = 200000
$ws.Range("G45").Value = 700000
$ws.Range("H45").Value = 700000
# Row 46
$ws.Range("D46").Value = 0
$ws.Range("E46").Value = 0
$ws.Range("F46").Value = 0
$ws.Range("G46").Value = 0
$ws.Range("H46").Value = 0
# Row 47
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("F47").Value = 471104
$ws.Range("G47").Value = 0
$ws.Range("H47").Value = 0
# Row 48
$ws.Range("D48").Value = 0
$ws.Range("E48").Value = 0
$ws.Range("F48").Value = 0
$ws.Range("G48").Value = 0
$ws.Range("H48").Value = 0
# Row 49
$ws.Range("D49").Value = 0
$ws.Range("E49").Value = 0
$ws.Range("F49").Value = 0
$ws.Range("G49").Value = 0
$ws.Range("H49").Value = 0
# Row 50
$ws.Range("D50").Value = 20000
$ws.Range("E50").Value = 20000
$ws.Range("F50").Value = 20000
$ws.Range("G50").Value = 70000
$ws.Range("H50").Value = 70000
# Row 51
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("F51").Value = 0
$ws.Range("G51").Value = 0
$ws.Range("H51").Value = 0
# Row 52
$ws.Range("D52").Value = "-"
$ws.Range("E52").Value = "-"
$ws.Range("F52").Value = "-"
$ws.Range("G52").Value = "-"
$ws.Range("H52").Value = "-"
# Row 53
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("F53").Value = 0
$ws.Range("G53").Value = 0
$ws.Range("H53").Value = 0
# Row 54
$ws.Range("D54").Value = "-"
$ws.Range("E54").Value = "-"
$ws.Range("F54").Value = "-"
$ws.Range("G54").Value = "-"
$ws.Range("H54").Value = "-"
# Row 55
$ws.Range("D55").Value = 0
$ws.Range("E55").Value = 0
$ws.Range("F55").Value = 0
$ws.Range("G55").Value = 0
$ws.Range("H55").Value = 0
# Row 56
$ws.Range("D56").Value = 768952
$ws.Range("E56").Value = 558626
$ws.Range("F56").Value = 681863
$ws.Range("G56").Value = 2351503
$ws.Range("H56").Value = 3522799
# Row 57
$ws.Range("D57").Value = 988952
$ws.Range("E57").Value = 778626
$ws.Range("F57").Value = 1372967
$ws.Range("G57").Value = 3121503
$ws.Range("H57").Value = 4292799
# Row 58
$ws.Range("D58").Value = 1369101
$ws.Range("E58").Value = 1534682
$ws.Range("F58").Value = 2967376
$ws.Range("G58").Value = 6593034
$ws.Range("H58").Value = 10526526